# Stroop replication trial types:
# - Rename header "corrAns" -> "correctAnswer"
# - Collapse "cong"/"incong" shorthand values into "congruent"/"incongruent"
# - Update the active selection from C6 to B5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename corrAns column header to correctAnswer
$ws.Range("C1").Value2 = "correctAnswer"

# Data rows: expand the abbreviated congruency values
$ws.Range("D2").Value2 = "congruent"
$ws.Range("D3").Value2 = "incongruent"
$ws.Range("D4").Value2 = "incongruent"
$ws.Range("D5").Value2 = "congruent"

# Update selected cell to B5
$ws.Range("B5").Select()
